$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    79  = "first day of the month"
    80  = "second day of the month"
    81  = "third day of the month"
    82  = "fourth day of the month"
    83  = "fifth day of the month"
    84  = "sixth day of the month"
    85  = "seventh day of the month"
    86  = "eighth day of the month"
    87  = "ninth day of the month"
    88  = "tenth day of the month"
    89  = "eleventh day of the month"
    90  = "twelth day of the month"
    91  = "thirteenth day of the month"
    92  = "fourteenth day of the month"
    93  = "fifteenth day of the month"
    94  = "sixteenth day of the month"
    95  = "seventeenth day of the month"
    96  = "eighteenth day of the month"
    97  = "nineteenth day of the month"
    98  = "twentieth day of the month"
    99  = "twenty-first day of the month"
    100 = "twenty-second day of the month"
    101 = "twenty-third day of the month"
    102 = "twenty-fourth day of the month"
    103 = "twenty-fifth day of the month"
    104 = "twenty-sixth day of the month"
    105 = "twenty-seventh day of the month"
    106 = "twenty-eighth day of the month"
    107 = "twenty-ninth day of the month"
    108 = "thirtieth day of the month"
    109 = "thiry-first day of the month"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
